$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D:D").Insert(-4161, 1)
Write-Host "done"
